$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save off the current (pre-edit) row 80 values - this data will be moved
# down to the newly inserted row 81.
$oldRow80 = @(
    $ws.Cells.Item(80, 1).Value2,
    $ws.Cells.Item(80, 2).Value2,
    $ws.Cells.Item(80, 3).Value2,
    $ws.Cells.Item(80, 4).Value2,
    $ws.Cells.Item(80, 5).Value2,
    $ws.Cells.Item(80, 6).Value2,
    $ws.Cells.Item(80, 7).Value2,
    $ws.Cells.Item(80, 8).Value2,
    $ws.Cells.Item(80, 9).Value2,
    $ws.Cells.Item(80, 10).Value2,
    $ws.Cells.Item(80, 11).Value2,
    $ws.Cells.Item(80, 12).Value2,
    $ws.Cells.Item(80, 13).Value2,
    $ws.Cells.Item(80, 14).Value2,
    $ws.Cells.Item(80, 15).Value2,
    $ws.Cells.Item(80, 16).Value2,
    $ws.Cells.Item(80, 17).Value2,
    $ws.Cells.Item(80, 18).Value2
)

# Insert a new blank row at 81, shifting the old row 81 (and below) down to 82.
$ws.Rows.Item(81).Insert(-4121)

# Populate the newly inserted row 81 with the data that used to live in row 80.
$ws.Cells.Item(81, 1).Value2 = $oldRow80[0]
$ws.Cells.Item(81, 2).Value2 = $oldRow80[1]
$ws.Cells.Item(81, 3).Value2 = $oldRow80[2]
$ws.Cells.Item(81, 4).Value2 = $oldRow80[3]
$ws.Cells.Item(81, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(81, 5).Value2 = $oldRow80[4]
$ws.Cells.Item(81, 6).Value2 = $oldRow80[5]
$ws.Cells.Item(81, 7).Value2 = $oldRow80[6]
$ws.Cells.Item(81, 8).Value2 = $oldRow80[7]
$ws.Cells.Item(81, 9).Value2 = $oldRow80[8]
$ws.Cells.Item(81, 10).Value2 = $oldRow80[9]
$ws.Cells.Item(81, 11).Value2 = $oldRow80[10]
$ws.Cells.Item(81, 12).Value2 = $oldRow80[11]
$ws.Cells.Item(81, 13).Value2 = $oldRow80[12]
$ws.Cells.Item(81, 14).Value2 = $oldRow80[13]
$ws.Cells.Item(81, 15).Value2 = $oldRow80[14]
$ws.Cells.Item(81, 16).Value2 = $oldRow80[15]
$ws.Cells.Item(81, 17).Value2 = $oldRow80[16]
$ws.Cells.Item(81, 18).Value2 = $oldRow80[17]

# Update row 80 with the new weekly price entry.
$ws.Cells.Item(80, 4).Value2 = 44595
$ws.Cells.Item(80, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(80, 10).Value2 = 80
$ws.Cells.Item(80, 11).Value2 = 23000
$ws.Cells.Item(80, 12).Value2 = 24000
$ws.Cells.Item(80, 13).Value2 = 23500
$ws.Cells.Item(80, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(80, 16).Value2 = 940
